# Apply updated dSF (column F) values to the data rows.
# The workbook's F column ("dSF") values were repulled/recalculated; this
# script writes the new values for the affected rows, leaving all other
# cells untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newF = @{
    2  = -1
    3  = -2
    4  = -1
    5  = 3
    6  = 3
    7  = -4
    8  = -5
    9  = 1
    10 = -2
    12 = 4
    14 = 3
    15 = -4
    16 = -2
    17 = -3
    19 = 1
    20 = 11
    22 = -1
    23 = -2
    24 = -4
    25 = 2
    26 = -1
    27 = 1
    28 = 3
    29 = -2
    30 = -3
    31 = -2
    33 = -2
    34 = 2
    35 = -2
}

foreach ($row in $newF.Keys) {
    $ws.Range("F$row").Value = $newF[$row]
}
